$d = $word.ActiveDocument

function AppendPara([string]$style) {
    $r = $d.Content
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $p = $d.Paragraphs.Last
    $p.Style = $style
    return $p
}

function AppendRun([string]$text) {
    $p = $d.Paragraphs.Last
    $ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
    $ip.InsertAfter($text)
}

function AppendItalicPara() {
    $r = $d.Content
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $p = $d.Paragraphs.Last
    $p.Style = "Normal"
    $p.Range.Italic = 1
    return $p
}

function AppendItalicRun([string]$text) {
    $p = $d.Paragraphs.Last
    $ip = $d.Range($p.Range.End - 1, $p.Range.End - 1)
    $ip.InsertAfter($text)
    $ip.Italic = 1
}

# --- 6 blank (Normal) paragraphs -----------------------------------------
AppendPara "Normal" | Out-Null
AppendPara "Normal" | Out-Null
AppendPara "Normal" | Out-Null
AppendPara "Normal" | Out-Null
AppendPara "Normal" | Out-Null
AppendPara "Normal" | Out-Null

# --- 2 empty "Titre1" paragraphs -----------------------------------------
AppendPara "Titre1" | Out-Null
AppendPara "Titre1" | Out-Null

# --- "Titre1" paragraph holding the report title --------------------------
AppendPara "Titre1" | Out-Null
AppendRun "Rapport d'activités"

# --- blank paragraph --------------------------------------------------------
AppendPara "Normal" | Out-Null

# --- "Titre2" paragraph holding the author's name --------------------------
AppendPara "Titre2" | Out-Null
AppendRun "Avram Iulian"

# --- blank paragraph --------------------------------------------------------
AppendPara "Normal" | Out-Null

# --- first italic quote paragraph -------------------------------------------
AppendItalicPara | Out-Null
AppendItalicRun '"'
AppendItalicRun "En plus de la participation "
AppendItalicRun "à"
AppendItalicRun " tous les cours, j'ai fais le diagramme des UC (60 minutes) et la partie ""Fonctionnalités"" du rapport "
AppendItalicRun "(60 minutes) et ma partie pour le rapport d'activités(5 minutes)."

# --- second italic quote paragraph ------------------------------------------
AppendItalicPara | Out-Null
AppendItalicRun "En classe, j'ai relu l'appel d'offre et avec l'ensemble du groupe on a choisi un nom pour notre application de plus durant les cours on a tous donné des conseils à Pawel pour améliorer les "
AppendItalicRun "IHMs"
AppendItalicRun "."
AppendItalicRun '"'
